# Financials update - yearly figures revised (commit: "Doing Updates for Financials")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total Revenue (row 8/9)
$ws.Range("G8").Value = 100
$ws.Range("G9").Value = 100

# Operating expenses / income section (rows 18-35)
$ws.Range("F18").Value = -700

$ws.Range("F21").Value = -700
$ws.Range("J21").Value = -300

$ws.Range("J22").Value = 400

$ws.Range("F23").Value = -700
$ws.Range("I23").Value = -900

$ws.Range("F26").Value = -700
$ws.Range("I26").Value = -900

$ws.Range("F27").Value = -700
$ws.Range("I27").Value = -900

$ws.Range("F33").Value = -700
$ws.Range("I33").Value = -900

$ws.Range("F35").Value = -700
$ws.Range("I35").Value = -900

# Balance sheet section (rows 57-76)
$ws.Range("F57").Value = 700
$ws.Range("G57").Value = 700
$ws.Range("H57").Value = 700

$ws.Range("D58").Value = 3800
$ws.Range("E58").Value = 3400
$ws.Range("F58").Value = 3000
$ws.Range("G58").Value = 2700
$ws.Range("H58").Value = 2400
$ws.Range("J58").Value = 1900

$ws.Range("D60").Value = 4700
$ws.Range("E60").Value = 4300
$ws.Range("F60").Value = 3900
$ws.Range("G60").Value = 3500
$ws.Range("H60").Value = 3200
$ws.Range("I60").Value = 2900
$ws.Range("J60").Value = 2500

$ws.Range("D61").Value = 4200
$ws.Range("E61").Value = 3800
$ws.Range("F61").Value = 3400
$ws.Range("G61").Value = 3000
$ws.Range("H61").Value = 2700
$ws.Range("I61").Value = 2400
$ws.Range("J61").Value = 2100

$ws.Range("D66").Value = 8900
$ws.Range("E66").Value = 8000
$ws.Range("F66").Value = 7300
$ws.Range("G66").Value = 6500
$ws.Range("H66").Value = 5900
$ws.Range("I66").Value = 5300
$ws.Range("J66").Value = 4600

$ws.Range("D72").Value = -57700
$ws.Range("E72").Value = -56800
$ws.Range("F72").Value = -56000
$ws.Range("G72").Value = -55300
$ws.Range("H72").Value = -54600
$ws.Range("I72").Value = -53900
$ws.Range("J72").Value = -53000

$ws.Range("D76").Value = -8900
$ws.Range("E76").Value = -8000
$ws.Range("F76").Value = -7300
$ws.Range("G76").Value = -6500
$ws.Range("H76").Value = -5800
$ws.Range("I76").Value = -5100
$ws.Range("J76").Value = -4300

# Cash flow statement section (rows 81-100)
$ws.Range("F81").Value = -700
$ws.Range("I81").Value = -900

# Capital Expenditures row: was "NA" text for D:H, now numeric 0; J goes from -300 to 0
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0

$ws.Range("J100").Value = 300
